# Apply the edits described in the commit "Error Calculations and Plots"
# to the missing_data worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Remove the "RM 232" row (originally row 26) and the "SC 92" row
#    (originally row 28, which becomes row 27 once row 26 is removed).
#    Deleting whole rows shifts everything below upward, matching the
#    row renumbering seen between the before/after sheets (35 -> 33
#    total rows).
# ---------------------------------------------------------------------
$ws.Rows.Item(26).Delete()   # removes "RM 232" (row 26); "SC 92" is now row 27
$ws.Rows.Item(27).Delete()   # removes "SC 92" (now at row 27)

# ---------------------------------------------------------------------
# 2) Scattered single-cell value changes across rows 2-24 (columns D/F)
# ---------------------------------------------------------------------
$ws.Cells.Item(2, 6).ClearContents()          # F2: 18.03 -> (blank)

$ws.Cells.Item(5, 6).Value = 17.66            # F5: (blank) -> 17.66

$ws.Cells.Item(6, 4).Value = -14.2            # D6: (blank) -> -14.2
$ws.Cells.Item(6, 6).Value = 16.43            # F6: (blank) -> 16.43

$ws.Cells.Item(8, 4).ClearContents()          # D8: -13.9 -> (blank)

$ws.Cells.Item(9, 6).ClearContents()          # F9: 17.26 -> (blank)

$ws.Cells.Item(10, 6).ClearContents()         # F10: 16.43 -> (blank)

$ws.Cells.Item(12, 4).Value = -14.1           # D12: (blank) -> -14.1

$ws.Cells.Item(14, 4).ClearContents()         # D14: -13.1 -> (blank)

$ws.Cells.Item(17, 4).Value = -14.7           # D17: (blank) -> -14.7

$ws.Cells.Item(18, 4).Value = -15.2           # D18: (blank) -> -15.2

$ws.Cells.Item(19, 4).ClearContents()         # D19: -15.5 -> (blank)

$ws.Cells.Item(20, 4).ClearContents()         # D20: -14 -> (blank)

$ws.Cells.Item(23, 4).Value = -13.9           # D23: (blank) -> -13.9

$ws.Cells.Item(24, 6).Value = 16.78           # F24: (blank) -> 16.78

# ---------------------------------------------------------------------
# 3) Fix-ups on the shifted "SC *" rows (now rows 26-33 after the two
#    row deletions above).
# ---------------------------------------------------------------------
$ws.Cells.Item(27, 2).Value = -20.4           # SC 101, B27: (blank) -> -20.4
$ws.Cells.Item(27, 4).ClearContents()         # SC 101, D27: -14.6 -> (blank)

$ws.Cells.Item(28, 2).ClearContents()         # SC 105, B28: -19.6 -> (blank)
$ws.Cells.Item(28, 6).ClearContents()         # SC 105, F28: 17.44 -> (blank)

$ws.Cells.Item(29, 2).ClearContents()         # SC 119, B29: -19.5 -> (blank)

$ws.Cells.Item(30, 2).Value = -19.7           # SC 120, B30: (blank) -> -19.7
$ws.Cells.Item(30, 6).Value = 16.89           # SC 120, F30: (blank) -> 16.89

$ws.Cells.Item(32, 2).ClearContents()         # SC 193, B32: -19.9 -> (blank)
